$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension/measure metadata after re-processing with curated dimensions.
$ws.Range("A2").Value = "iaest-measure:tipo-de-hogar"
$ws.Range("C2").Value = "sdmx-dimension:refArea"

$ws.Range("A3").Value = "medida"
$ws.Range("C3").Value = "dim"

$ws.Range("A4").Value = "xsd:int"
$ws.Range("C4").Value = "URI-Municipio"

# Remove the now-obsolete mapping-file row.
$ws.Range("A5").EntireRow.Delete()
